$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LuckyNumber")
$ws.Name = "Giải thưởng"
